# Implemented getting kafka relations.
# The "classFields" sheet's field listing order changed (re-generated from a
# reflection-based structure scan), so the Field Name / Field Type columns
# for each row need to be rewritten to reflect the new ordering while the
# Class Name / Field Modifier columns stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Row => @(FieldName, FieldType)
$fieldData = @{
    2  = @("name", "java.lang.String")
    3  = @("serialVersionUID", "long")
    4  = @("type", "java.lang.Integer")
    5  = @("value", "java.lang.String")
    6  = @("productAttributeId", "java.lang.Long")
    7  = @("id", "java.lang.Long")
    8  = @("esProductService", "com.macro.mall.search.service.EsProductService")
    9  = @("id", "java.lang.Long")
    10 = @("sale", "java.lang.Integer")
    11 = @("brandId", "java.lang.Long")
    12 = @("newStatus", "java.lang.Integer")
    13 = @("stock", "java.lang.Integer")
    14 = @("sort", "java.lang.Integer")
    15 = @("attrValueList", "java.util.List")
    16 = @("promotionType", "java.lang.Integer")
    17 = @("productSn", "java.lang.String")
    18 = @("keywords", "java.lang.String")
    19 = @("productCategoryName", "java.lang.String")
    20 = @("recommandStatus", "java.lang.Integer")
    21 = @("productCategoryId", "java.lang.Long")
    22 = @("serialVersionUID", "long")
    23 = @("brandName", "java.lang.String")
    24 = @("name", "java.lang.String")
    25 = @("pic", "java.lang.String")
    26 = @("subTitle", "java.lang.String")
    27 = @("price", "java.math.BigDecimal")
    28 = @("elasticsearchRestTemplate", "org.springframework.data.elasticsearch.core.ElasticsearchRestTemplate")
    29 = @("productRepository", "com.macro.mall.search.repository.EsProductRepository")
    30 = @("LOGGER", "org.slf4j.Logger")
    31 = @("productDao", "com.macro.mall.search.dao.EsProductDao")
    32 = @("attrValues", "java.util.List")
    33 = @("attrName", "java.lang.String")
    34 = @("attrId", "java.lang.Long")
    35 = @("productCategoryNames", "java.util.List")
    36 = @("productAttrs", "java.util.List")
    37 = @("brandNames", "java.util.List")
}

foreach ($row in $fieldData.Keys) {
    $values = $fieldData[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 4).Value = $values[1]
}
